$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 (RM 232) entirely - rows below shift up
$ws.Rows.Item(26).Delete()

# After the above deletion, the row that held "SC 92" is now row 27 - delete it too
$ws.Rows.Item(27).Delete()

# Fill in previously-missing values in column E (now at their shifted positions)
$ws.Range("E26").Value = -5        # row for "SC 5"
$ws.Range("E27").Value = $null     # row for "SC 101" -> now empty
$ws.Range("E33").Value = -10.7     # row for "SC 232"
